$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forms")

# --- Update the three "registrant email" cells (B2:B4) to the new address ---
$ws.Range("R6").Value2 = "'0"
$ws.Range("B2").Value2 = "qatesting.lotuswave+1@gmail.com"
$ws.Range("B3").Value2 = "qatesting.lotuswave+1@gmail.com"
$ws.Range("B4").Value2 = "qatesting.lotuswave+1@gmail.com"

# --- Add the new "Quantity" column header in R1, matching the formatting of the neighboring header cells ---
$ws.Range("R1").Value2 = "Quantity"
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)

# --- Add the new "Product Qunatity" row (row 6) ---
$ws.Range("A6").Value2 = "Product Qunatity"

# --- Rebuild the hyperlinks: delete all, then re-add them in the desired order / with the desired targets ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:qatesting.lotuswave@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:qatesting.lotuswave@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H4"), "mailto:qatesting.lotuswave@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:qatesting.lotuswave+1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:qatesting.lotuswave+1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:qatesting.lotuswave+1@gmail.com")

# --- Update the active selection shown when the workbook is reopened ---
$ws.Activate()
$ws.Range("E12").Select()
